$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 632.75757
$ws.Range("I6").Value = 106.703705
$ws.Range("J6").Value = 3000
$ws.Range("K6").Value = 320.111115
$ws.Range("L6").Value = 9000
$ws.Range("M6").Value = -208.111115
$ws.Range("N6").Value = -9224
$ws.Range("H17").Value = 230.76
$ws.Range("J17").Value = 237.04167
$ws.Range("L17").Value = 711.12501
$ws.Range("N17").Value = -1047.12501
$ws.Range("H32").Value = 669.8946999999999
$ws.Range("I32").Value = 500.25
$ws.Range("J32").Value = 715.13336
$ws.Range("K32").Value = 500.25
$ws.Range("L32").Value = 715.13336
$ws.Range("M32").Value = -174.25
$ws.Range("N32").Value = -1367.13336
$ws.Range("H40").Value = 1641.6666
$ws.Range("I40").Value = 1490
$ws.Range("J40").Value = 1750
$ws.Range("K40").Value = 1490
$ws.Range("L40").Value = 1750
$ws.Range("M40").Value = -1315
$ws.Range("N40").Value = -2100
$ws.Range("H48").Value = 1839
$ws.Range("J48").Value = 2250
$ws.Range("L48").Value = 6750
$ws.Range("N48").Value = -7334
$ws.Range("H56").Value = 1839
$ws.Range("J56").Value = 2250
$ws.Range("L56").Value = 6750
$ws.Range("N56").Value = -7818
$ws.Range("H64").Value = 3215.8572
$ws.Range("I64").Value = 2904.2
$ws.Range("J64").Value = 3995
$ws.Range("K64").Value = 2904.2
$ws.Range("L64").Value = 3995
$ws.Range("M64").Value = -2656.2
$ws.Range("N64").Value = -4491
$ws.Range("H67").Value = 3215.8572
$ws.Range("I67").Value = 2904.2
$ws.Range("J67").Value = 3995
$ws.Range("K67").Value = 2904.2
$ws.Range("L67").Value = 3995
$ws.Range("M67").Value = -2046.2
$ws.Range("N67").Value = -5711
$ws.Range("H70").Value = 955.14813
$ws.Range("I70").Value = 750.5
$ws.Range("J70").Value = 1118.8667
$ws.Range("K70").Value = 2251.5
$ws.Range("L70").Value = 3356.6001
$ws.Range("M70").Value = -1981.5
$ws.Range("N70").Value = -3896.6001
$ws.Range("H73").Value = 955.14813
$ws.Range("I73").Value = 750.5
$ws.Range("J73").Value = 1118.8667
$ws.Range("K73").Value = 2251.5
$ws.Range("L73").Value = 3356.6001
$ws.Range("M73").Value = -1315.5
$ws.Range("N73").Value = -5228.6001
$ws.Range("H100").Value = 2712.5
$ws.Range("I100").Value = 2740
$ws.Range("J100").Value = 2666.6667
$ws.Range("K100").Value = 2740
$ws.Range("L100").Value = 2666.6667
$ws.Range("M100").Value = -2199
$ws.Range("N100").Value = -3748.6667
$ws.Range("H107").Value = 2307.2
$ws.Range("I107").Value = 2366.0688
$ws.Range("J107").Value = 600
$ws.Range("K107").Value = 2366.0688
$ws.Range("L107").Value = 600
$ws.Range("M107").Value = -446.0688
$ws.Range("N107").Value = -4440
$ws.Range("H137").Value = 1525.05
$ws.Range("I137").Value = 949.75
$ws.Range("J137").Value = 1668.875
$ws.Range("K137").Value = 2849.25
$ws.Range("L137").Value = 5006.625
$ws.Range("M137").Value = -299.25
$ws.Range("N137").Value = -10106.625
$ws.Range("H138").Value = 5193.622
$ws.Range("J138").Value = 4971.4736
$ws.Range("L138").Value = 14914.4208
$ws.Range("N138").Value = -25194.4208

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 500
$ws.Range("I11").Value = 500
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 500
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()
$ws.Range("M11").Value = -356
$ws.Range("H32").Value = 501353.28
$ws.Range("I32").Value = 663903.1
$ws.Range("J32").Value = 13703.786
$ws.Range("K32").Value = 663903.1
$ws.Range("L32").Value = 13703.786
$ws.Range("M32").Value = -663616.1
$ws.Range("N32").Value = -14277.786
$ws.Range("H61").Value = 2524.1785
$ws.Range("I61").Value = 2054.6191
$ws.Range("J61").Value = 3932.8572
$ws.Range("K61").Value = 2054.6191
$ws.Range("L61").Value = 3932.8572
$ws.Range("M61").Value = -1842.6191
$ws.Range("N61").Value = -4356.8572
$ws.Range("H74").Value = 1674.8
$ws.Range("I74").Value = 1735.3
$ws.Range("J74").Value = 1614.3
$ws.Range("K74").Value = 1735.3
$ws.Range("L74").Value = 1614.3
$ws.Range("M74").Value = -861.3
$ws.Range("N74").Value = -3362.3
$ws.Range("H77").Value = 1674.8
$ws.Range("I77").Value = 1735.3
$ws.Range("J77").Value = 1614.3
$ws.Range("K77").Value = 8676.5
$ws.Range("L77").Value = 8071.5
$ws.Range("M77").Value = -4308.5
$ws.Range("N77").Value = -16807.5
$ws.Range("H132").Value = 3440.1765
$ws.Range("I132").Value = 2629.025
$ws.Range("K132").Value = 7887.075000000001
$ws.Range("M132").Value = -5357.075000000001
$ws.Range("H134").Value = 71914.5
$ws.Range("J134").Value = 71914.5
$ws.Range("L134").Value = 71914.5
$ws.Range("N134").Value = -82054.5
$ws.Range("H136").Value = 2524.1785
$ws.Range("I136").Value = 2054.6191
$ws.Range("J136").Value = 3932.8572
$ws.Range("K136").Value = 6163.8573
$ws.Range("L136").Value = 11798.5716
$ws.Range("M136").Value = -3613.8573
$ws.Range("N136").Value = -16898.5716

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3683.3872
$ws.Range("I134").Value = 3996.4443
$ws.Range("J134").Value = 3249.923
$ws.Range("K134").Value = 11989.3329
$ws.Range("L134").Value = 9749.769
$ws.Range("M134").Value = -9454.332900000001
$ws.Range("N134").Value = -14819.769

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 677.75
$ws.Range("I16").Value = 605.5
$ws.Range("J16").Value = 750
$ws.Range("K16").Value = 605.5
$ws.Range("L16").Value = 750
$ws.Range("M16").Value = -318.5
$ws.Range("N16").Value = -1324
$ws.Range("H31").Value = 1514.2059
$ws.Range("I31").Value = 1639.3529
$ws.Range("J31").Value = 1389.0588
$ws.Range("K31").Value = 1639.3529
$ws.Range("L31").Value = 1389.0588
$ws.Range("M31").Value = -1344.3529
$ws.Range("N31").Value = -1979.0588
$ws.Range("H34").Value = 1514.2059
$ws.Range("I34").Value = 1639.3529
$ws.Range("J34").Value = 1389.0588
$ws.Range("K34").Value = 1639.3529
$ws.Range("L34").Value = 1389.0588
$ws.Range("M34").Value = -1437.3529
$ws.Range("N34").Value = -1793.0588
$ws.Range("H113").Value = 677.75
$ws.Range("I113").Value = 605.5
$ws.Range("J113").Value = 750
$ws.Range("K113").Value = 605.5
$ws.Range("L113").Value = 750
$ws.Range("M113").Value = 1564.5
$ws.Range("N113").Value = -5090
$ws.Range("H122").Value = 1427.32
$ws.Range("I122").Value = 1457.6428
$ws.Range("J122").Value = 1388.7273
$ws.Range("K122").Value = 4372.928400000001
$ws.Range("L122").Value = 4166.1819
$ws.Range("M122").Value = -1922.928400000001
$ws.Range("N122").Value = -9066.1819
$ws.Range("H132").Value = 20836372
$ws.Range("I132").Value = 3162.2
$ws.Range("K132").Value = 9486.599999999999
$ws.Range("M132").Value = -6956.599999999999
$ws.Range("H134").Value = 1575.1428
$ws.Range("I134").Value = 1105.2
$ws.Range("K134").Value = 3315.6
$ws.Range("M134").Value = -780.6000000000004

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 632
$ws.Range("H23").Value = 76923150
$ws.Range("I23").Value = 62.2
$ws.Range("J23").Value = 125000080
$ws.Range("K23").Value = 186.6
$ws.Range("L23").Value = 375000240
$ws.Range("M23").Value = 48.39999999999998
$ws.Range("N23").Value = -375000710
$ws.Range("H68").Value = 848.7143
$ws.Range("I68").Value = 609.8946999999999
$ws.Range("J68").Value = 971.3514
$ws.Range("K68").Value = 1829.6841
$ws.Range("L68").Value = 2914.0542
$ws.Range("M68").Value = -1018.6841
$ws.Range("N68").Value = -4536.0542
$ws.Range("H71").Value = 848.7143
$ws.Range("I71").Value = 609.8946999999999
$ws.Range("J71").Value = 971.3514
$ws.Range("K71").Value = 5489.052299999999
$ws.Range("L71").Value = 8742.1626
$ws.Range("M71").Value = -1433.052299999999
$ws.Range("N71").Value = -16854.1626
$ws.Range("H86").Value = 1670
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 1670
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 5010
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -7382
$ws.Range("H89").Value = 1670
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 1670
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 15030
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -26886
$ws.Range("H112").Value = 5390.6294
$ws.Range("I112").Value = 3309
$ws.Range("J112").Value = 5650.8335
$ws.Range("K112").Value = 9927
$ws.Range("L112").Value = 16952.5005
$ws.Range("M112").Value = -8819
$ws.Range("N112").Value = -19168.5005
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("M128").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3646.4285
$ws.Range("I132").Value = 3606.5
$ws.Range("J132").Value = 3699.6667
$ws.Range("K132").Value = 10819.5
$ws.Range("L132").Value = 11099.0001
$ws.Range("M132").Value = -8289.5
$ws.Range("N132").Value = -16159.0001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H94").Value = 62500
$ws.Range("J94").Value = 62500
$ws.Range("L94").Value = 62500
$ws.Range("N94").Value = -63852
$ws.Range("H132").Value = 3820.32
$ws.Range("I132").Value = 3107.8572
$ws.Range("K132").Value = 9323.571599999999
$ws.Range("M132").Value = -6793.571599999999
$ws.Range("H136").Value = 1541
$ws.Range("I136").Value = 1541
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4623
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2073
$ws.Range("N136").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 15157293
$ws.Range("I132").Value = 6943.875
$ws.Range("J132").Value = 55558224
$ws.Range("K132").Value = 20831.625
$ws.Range("L132").Value = 166674672
$ws.Range("M132").Value = -18301.625
$ws.Range("N132").Value = -166679732
